$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 0.234
$ws.Range("M2").Value = 11.33322749852823
$ws.Range("E3").Value = 0.187
$ws.Range("M3").Value = 8.050491762716756
$ws.Range("E4").Value = 0.374
$ws.Range("M4").Value = 12.83644853621241
$ws.Range("E5").Value = 0.358
$ws.Range("M5").Value = 4.503817435818375
$ws.Range("E6").Value = 1.148
$ws.Range("M6").Value = 8.669311379419982
$ws.Range("E7").Value = 1.095
$ws.Range("M7").Value = 12.00567069656942
$ws.Range("E8").Value = 1.224
$ws.Range("M8").Value = 21.00548105532981
$ws.Range("E9").Value = 2.215
$ws.Range("M9").Value = 15.07025496826344
$ws.Range("E10").Value = 4.615
$ws.Range("M10").Value = 10.69303885368356
$ws.Range("E11").Value = 7.146999999999999
$ws.Range("M11").Value = 30.2653391232169
$ws.Range("E12").Value = 21.196
$ws.Range("M12").Value = 23.02755400127739
$ws.Range("E13").Value = 9.736
$ws.Range("M13").Value = 20.77724461659307
$ws.Range("E14").Value = 13.939
$ws.Range("M14").Value = 21.202664330512
$ws.Range("E15").Value = 4.598
$ws.Range("M15").Value = 4.018895980584129
$ws.Range("E16").Value = 171.666
$ws.Range("M16").Value = 57.21285957363571
$ws.Range("E17").Value = 76.313
$ws.Range("M17").Value = 67.57404752132992
$ws.Range("E18").Value = 78.456
$ws.Range("M18").Value = 42.51151224517732
$ws.Range("E19").Value = 17.721
$ws.Range("M19").Value = 4.021950074385358
$ws.Range("E20").Value = 2807.385
$ws.Range("M20").Value = 36.05133808510136
$ws.Range("E21").Value = 1904.577
$ws.Range("M21").Value = 48.9838631096483
